$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing rows 135-141 ---

# Row 135: title date text changed
$ws.Range("A135").Value = "日期：2018.11.07 第十周周三"

# Row 137: Wang Weifeng completion % filled in
$ws.Range("C137").Value = 1

# Row 138: Chen Shengyun completion % filled in
$ws.Range("C138").Value = 0.95

# Row 139: Lin Weicheng completion % filled in
$ws.Range("C139").Value = 1

# Row 140: Wu Shuaichen remark added
$ws.Range("D140").Value = "优化情况不计入任务完成情况"

# Row 141: Li Haiyang completion % filled in
$ws.Range("C141").Value = 1

# --- Append a brand new weekly block (rows 145-153), cloned from the
#     135-143 block so borders/fonts/number-formats match exactly.
#     Pre-merge the destination title/footer ranges on their default
#     (unformatted) styles first, THEN paste the source formats on top -
#     this keeps the merged ranges while avoiding the host from minting
#     brand-new style entries for the merge. ---

$ws.Range("A145:D145").Merge()
$ws.Range("A152:D153").Merge()

$ws.Range("A135:D143").Copy()
$ws.Range("A145:D153").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Title
$ws.Range("A145").Value = "日期：2018.11.08 第十周周四"

# Header row 146 (组员/计划内容/完成情况/备注)
$ws.Range("A146").Value = "组员"
$ws.Range("B146").Value = "计划内容"
$ws.Range("C146").Value = "完成情况"
$ws.Range("D146").Value = "备注"

# Row 147 - Wang Weifeng
$ws.Range("A147").Value = "王伟锋"
$ws.Range("B147").Value = "基本功能已完成，现优后台服务系统"
$ws.Range("D147").Value = "优化情况不计入任务完成情况"

# Row 148 - Chen Shengyun
$ws.Range("A148").Value = "陈升云"
$ws.Range("B148").Value = "完成个人资料的查看，修改等"

# Row 149 - Lin Weicheng
$ws.Range("A149").Value = "林玮成"
$ws.Range("B149").Value = "上一阶段的android开发完成单元的测试文档开始编写"

# Row 150 - Wu Shuaichen
$ws.Range("A150").Value = "吴帅辰"
$ws.Range("B150").Value = "基本功能已完成，现优化管理员系统"
$ws.Range("D150").Value = "优化情况不计入任务完成情况"

# Row 151 - Li Haiyang
$ws.Range("A151").Value = "李海洋"
$ws.Range("B151").Value = "完成群简介的获取"

# Row 152 - summary label (carried over by the copy, but set explicitly
# to be safe)
$ws.Range("A152").Value = "总结："

# Match the author's final cursor position in the saved workbook
$ws.Range("E148").Select()

Write-Host "edit applied"
